# Munición Incendiaria DDBB v0.12
#
# The edit inserts one new card row at row 11 of "Hoja1" (pushing every
# row from 11 downward down by one, 558 -> 559 total rows), for the new
# card "Municion Incendiaria" / "da +2/+0 a los aliados".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a fresh blank row at row 11; everything currently at 11.. shifts to 12..
$ws.Rows("11:11").Insert()

# The freshly inserted row has no formatting of its own yet - pull the
# cell formatting (borders/fill/font/alignment/number format) for A:G from
# the row that now holds the old row-11 card (row 12), and H's formatting
# from a row that uses the wrapped "ability text" style (row 13, style 9),
# matching the layout of every other card row in the table.
$ws.Range("A12:G12").Copy() | Out-Null
$ws.Range("A11:G11").PasteSpecial(-4122) | Out-Null

$ws.Range("H13").Copy() | Out-Null
$ws.Range("H11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Match the row height used by the other single-line card rows.
$ws.Rows("11:11").RowHeight = 15.75

# Fill in the new card's data.
$ws.Range("A11").Value = "Municion Incendiaria"
$ws.Range("B11").Value = $ws.Range("B12").Value()
$ws.Range("C11").Value = $ws.Range("C12").Value()
$ws.Range("D11").Value = $ws.Range("D12").Value()
$ws.Range("E11").Value = $ws.Range("E12").Value()
$ws.Range("F11").Value = $ws.Range("F12").Value()
$ws.Range("G11").Value = 3
$ws.Range("H11").Value = "da +2/+0 a los aliados"

# Leave selection where the author left it.
$ws.Range("H11").Select()
